$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adm"
$ws.Range("C2").Value = "Calcr"
$ws.Range("D2").Value = "Neutro"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.455024
$ws.Range("H2").Value = 24.910048
$ws.Range("I2").Value = 0.2137786973425757
$ws.Range("J2").Value = 0.1547465300518145
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.008382333333333334
$ws.Range("N2").Value = 0.025147
$ws.Range("O2").Value = 0.00230362653310207
$ws.Range("P2").Value = 0.003451464357217434
$ws.Range("Q2").Value = 0.1044021628426667
$ws.Range("R2").Value = 0.6264129770560001
$ws.Range("S2").Value = 0.0004924662794103544
$ws.Range("T2").Value = 0.0005341021328769144

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adm"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.455024
$ws.Range("H3").Value = 24.910048
$ws.Range("I3").Value = 0.2137786973425757
$ws.Range("J3").Value = 0.1547465300518145
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.630373
$ws.Range("N3").Value = 7.260746
$ws.Range("O3").Value = 0.9976963734668979
$ws.Range("P3").Value = 0.9965485356427827
$ws.Range("Q3").Value = 45.216382843952
$ws.Range("R3").Value = 180.865531375808
$ws.Range("S3").Value = 0.2132862310631654
$ws.Range("T3").Value = 0.1542124279189377

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Adm"
$ws.Range("C4").Value = "Calcr"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 32.753573
$ws.Range("H4").Value = 98.260719
$ws.Range("I4").Value = 0.5621840768235339
$ws.Range("J4").Value = 0.6104165397692691
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008382333333333334
$ws.Range("N4").Value = 0.025147
$ws.Range("O4").Value = 0.00230362653310207
$ws.Range("P4").Value = 0.003451464357217434
$ws.Range("Q4").Value = 0.2745513667436666
$ws.Range("R4").Value = 2.470962300693
$ws.Range("S4").Value = 0.001295062155858185
$ws.Range("T4").Value = 0.00210683093006963

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Adm"
$ws.Range("C5").Value = "Calcr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 32.753573
$ws.Range("H5").Value = 98.260719
$ws.Range("I5").Value = 0.5621840768235339
$ws.Range("J5").Value = 0.6104165397692691
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.630373
$ws.Range("N5").Value = 7.260746
$ws.Range("O5").Value = 0.9976963734668979
$ws.Range("P5").Value = 0.9965485356427827
$ws.Range("Q5").Value = 118.907687072729
$ws.Range("R5").Value = 713.446122436374
$ws.Range("S5").Value = 0.5608890146676757
$ws.Range("T5").Value = 0.6083097088391995

# Row 6
$ws.Range("A6").Value = "Neutro"
$ws.Range("B6").Value = "Adm"
$ws.Range("C6").Value = "Calcr"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.69705966666667
$ws.Range("H6").Value = 35.091179
$ws.Range("I6").Value = 0.2007689570311854
$ws.Range("J6").Value = 0.2179938868715589
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.008382333333333334
$ws.Range("N6").Value = 0.025147
$ws.Range("O6").Value = 0.00230362653310207
$ws.Range("P6").Value = 0.003451464357217434
$ws.Range("Q6").Value = 0.09804865314588888
$ws.Range("R6").Value = 0.8824378783129999
$ws.Range("S6").Value = 0.000462496696440268
$ws.Range("T6").Value = 0.000752398130628475

# Row 7
$ws.Range("A7").Value = "Neutro"
$ws.Range("B7").Value = "Adm"
$ws.Range("C7").Value = "Calcr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.69705966666667
$ws.Range("H7").Value = 35.091179
$ws.Range("I7").Value = 0.2007689570311854
$ws.Range("J7").Value = 0.2179938868715589
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.630373
$ws.Range("N7").Value = 7.260746
$ws.Range("O7").Value = 0.9976963734668979
$ws.Range("P7").Value = 0.9965485356427827
$ws.Range("Q7").Value = 42.46468959325566
$ws.Range("R7").Value = 254.788137559534
$ws.Range("S7").Value = 0.2003064603347451
$ws.Range("T7").Value = 0.2172414887409305

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Adm"
$ws.Range("C8").Value = "Calcr"
$ws.Range("D8").Value = "Neutro"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.3556395
$ws.Range("H8").Value = 2.711279
$ws.Range("I8").Value = 0.02326826880270489
$ws.Range("J8").Value = 0.01684304330735749
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.008382333333333334
$ws.Range("N8").Value = 0.025147
$ws.Range("O8").Value = 0.00230362653310207
$ws.Range("P8").Value = 0.003451464357217434
$ws.Range("Q8").Value = 0.01136342216883333
$ws.Range("R8").Value = 0.068180533013
$ws.Range("S8").Value = 0.00005360140139326212
$ws.Range("T8").Value = 0.000058133163642414

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Adm"
$ws.Range("C9").Value = "Calcr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.3556395
$ws.Range("H9").Value = 2.711279
$ws.Range("I9").Value = 0.02326826880270489
$ws.Range("J9").Value = 0.01684304330735749
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.630373
$ws.Range("N9").Value = 7.260746
$ws.Range("O9").Value = 0.9976963734668979
$ws.Range("P9").Value = 0.9965485356427827
$ws.Range("Q9").Value = 4.9214770385335
$ws.Range("R9").Value = 19.685908154134
$ws.Range("S9").Value = 0.02321466740131163
$ws.Range("T9").Value = 0.01678491014371507
